$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold text-formatted numeric-looking
# strings (e.g. "65.553.76", "  +2.64%  "). Force the whole range to
# Text format first so assigning these values doesn't get silently
# auto-converted to numbers by Excel's smart-entry parsing; restore the
# original (default) style afterwards so cell formatting is unchanged.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "65.553.76"
$ws.Range("E2").Value = "  +2.64%  "
$ws.Range("D3").Value = "3.201.97"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "599.63"
$ws.Range("E5").Value = "  +1.87%  "
$ws.Range("D6").Value = "152.92"
$ws.Range("E6").Value = "  +4.59%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.199.08"
$ws.Range("E8").Value = "  +1.27%  "
$ws.Range("D9").Value = "0.544"
$ws.Range("E9").Value = "  +2.45%  "
$ws.Range("D10").Value = "0.167"
$ws.Range("E10").Value = "  +3.50%  "
$ws.Range("D11").Value = "6.11"
$ws.Range("E11").Value = "  +6.21%  "
$ws.Range("D12").Value = "0.472"
$ws.Range("E12").Value = "  +2.14%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "0.0000254"
$ws.Range("E13").Value = "  +2.47%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "39.42"
$ws.Range("E14").Value = "  +6.42%  "
$ws.Range("D15").Value = "3.739.34"
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "7.39"
$ws.Range("E17").Value = "  +3.64%  "
$ws.Range("D18").Value = "65.163.73"
$ws.Range("E18").Value = "  +2.24%  "
$ws.Range("D19").Value = "3.212.62"
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("D20").Value = "483.24"
$ws.Range("E20").Value = "  +3.78%  "
$ws.Range("D21").Value = "14.97"
$ws.Range("E21").Value = "  +4.47%  "
$ws.Range("D22").Value = "0.773"
$ws.Range("E22").Value = "  +5.57%  "
$ws.Range("D23").Value = "7.89"
$ws.Range("E23").Value = "  +5.25%  "
$ws.Range("D24").Value = "2.46"
$ws.Range("E24").Value = "  +12.00%  "
$ws.Range("D25").Value = "13.68"
$ws.Range("E25").Value = "  +4.67%  "
$ws.Range("D26").Value = "83.50"
$ws.Range("E26").Value = "  +2.47%  "
$ws.Range("E27").Value = "  +0.34%  "
$ws.Range("D28").Value = "9.85"
$ws.Range("E28").Value = "  +7.94%  "
$ws.Range("D29").Value = "2.79"
$ws.Range("E29").Value = "  +3.65%  "
$ws.Range("E30").Value = "  +3.34%  "
$ws.Range("D31").Value = "7.47"
$ws.Range("E31").Value = "  +6.36%  "
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("D33").Value = "0.122"
$ws.Range("E33").Value = "  +10.05%  "
$ws.Range("D34").Value = "28.67"
$ws.Range("E34").Value = "  +5.80%  "
$ws.Range("D35").Value = "0.0₃0900"
$ws.Range("E35").Value = "  +4.18%  "
$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D36").Value = "3.57"
$ws.Range("E36").Value = "  +5.52%  "
$ws.Range("B37").Value = "Mantle"
$ws.Range("C37").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D37").Value = "1.09"
$ws.Range("E37").Value = "  +4.25%  "
$ws.Range("D38").Value = "2.38"
$ws.Range("E38").Value = "  +2.89%  "
$ws.Range("D39").Value = "6.32"
$ws.Range("E39").Value = "  +4.75%  "
$ws.Range("D40").Value = "478.89"
$ws.Range("E40").Value = "  +8.00%  "
$ws.Range("D41").Value = "52.18"
$ws.Range("E41").Value = "  +3.16%  "
$ws.Range("D42").Value = "9.44"
$ws.Range("E42").Value = "  +8.19%  "
$ws.Range("D43").Value = "0.302"
$ws.Range("E43").Value = "  +10.01%  "
$ws.Range("D44").Value = "0.0385"
$ws.Range("E44").Value = "  +3.45%  "
$ws.Range("D45").Value = "2.944.46"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("E46").Value = "  +3.81%  "
$ws.Range("D47").Value = "38.92"
$ws.Range("E47").Value = "  +7.61%  "
$ws.Range("D48").Value = "131.63"
$ws.Range("E48").Value = "  +4.77%  "
$ws.Range("D49").Value = "2.33"
$ws.Range("E49").Value = "  +6.77%  "
$ws.Range("D50").Value = "25.40"
$ws.Range("E50").Value = "  +3.92%  "
$ws.Range("E51").Value = "  +0.00%  "

$ws.Range("D2:E51").Style = "Normal"
